$wb = $excel.ActiveWorkbook

# Delete the third sheet ("THIRD2")
$ws3 = $wb.Worksheets.Item("THIRD2")
$ws3.Delete()

# Rename remaining sheets
$ws1 = $wb.Worksheets.Item("FIRST2")
$ws1.Name = "FIRST"

$ws2 = $wb.Worksheets.Item("SECOND2")
$ws2.Name = "SECOND"

# Make the second sheet ("SECOND") the active/selected sheet
$ws2.Activate()
$ws2.Select()
